# Update "想去人数" (want-to-go counts) figures to the newly scraped values.
$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibitions)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 1882   # 南宁·AP动漫游戏嘉年华: 1876 -> 1882
$ws1.Range("F5").Value = 1570   # 南宁·布谷鸟动漫展4th: 1568 -> 1570
$ws1.Range("F8").Value = 557    # 南宁·AB动漫游戏嘉年华: 552 -> 557

# Sheet "全部类型" (all types) - same events, rows shifted by one due to an
# additional performance-type row present only in this aggregated sheet.
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 1882   # 南宁·AP动漫游戏嘉年华: 1876 -> 1882
$ws4.Range("F5").Value = 1570   # 南宁·布谷鸟动漫展4th: 1568 -> 1570
$ws4.Range("F9").Value = 557    # 南宁·AB动漫游戏嘉年华: 552 -> 557

$wb.Save()
